$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40,8).Value = 3918.5  # H40: 3969.5334 -> 3918.5
$ws.Cells.Item(40,9).Value = 3883.3333  # I40: 3645.8572 -> 3883.3333
$ws.Cells.Item(40,10).Value = 3953.6667  # J40: 4252.75 -> 3953.6667
$ws.Cells.Item(40,11).Value = 3883.3333  # K40: 3645.8572 -> 3883.3333
$ws.Cells.Item(40,12).Value = 3953.6667  # L40: 4252.75 -> 3953.6667
$ws.Cells.Item(40,13).Value = -3708.3333  # M40: -3470.8572 -> -3708.3333
$ws.Cells.Item(40,14).Value = -4303.6667  # N40: -4602.75 -> -4303.6667
$ws.Cells.Item(51,8).Value = 3989.6667  # H51: 5646 -> 3989.6667
$ws.Cells.Item(51,9).Value = 2000  # I51: 0 -> 2000
$ws.Cells.Item(51,10).Value = 4984.5  # J51: 5646 -> 4984.5
$ws.Cells.Item(51,11).Value = 2000  # K51: 0 -> 2000
$ws.Cells.Item(51,12).Value = 4984.5  # L51: 5646 -> 4984.5
$ws.Cells.Item(51,13).Value = -1516  # M51: None -> -1516
$ws.Cells.Item(51,14).Value = -5952.5  # N51: -6614 -> -5952.5
$ws.Cells.Item(76,8).Value = 7663.5884  # H76: 7423.5557 -> 7663.5884
$ws.Cells.Item(76,9).Value = 8365.666999999999  # I76: 8384.444 -> 8365.666999999999
$ws.Cells.Item(76,10).Value = 6873.75  # J76: 6462.6665 -> 6873.75
$ws.Cells.Item(76,11).Value = 8365.666999999999  # K76: 8384.444 -> 8365.666999999999
$ws.Cells.Item(76,12).Value = 6873.75  # L76: 6462.6665 -> 6873.75
$ws.Cells.Item(76,13).Value = -8050.666999999999  # M76: -8069.444 -> -8050.666999999999
$ws.Cells.Item(76,14).Value = -7503.75  # N76: -7092.6665 -> -7503.75
$ws.Cells.Item(79,8).Value = 7663.5884  # H79: 7423.5557 -> 7663.5884
$ws.Cells.Item(79,9).Value = 8365.666999999999  # I79: 8384.444 -> 8365.666999999999
$ws.Cells.Item(79,10).Value = 6873.75  # J79: 6462.6665 -> 6873.75
$ws.Cells.Item(79,11).Value = 8365.666999999999  # K79: 8384.444 -> 8365.666999999999
$ws.Cells.Item(79,12).Value = 6873.75  # L79: 6462.6665 -> 6873.75
$ws.Cells.Item(79,13).Value = -7273.666999999999  # M79: -7292.444 -> -7273.666999999999
$ws.Cells.Item(79,14).Value = -9057.75  # N79: -8646.666499999999 -> -9057.75
$ws.Cells.Item(98,8).Value = 90911980  # H98: 111113784 -> 90911980
$ws.Cells.Item(98,10).Value = 3933.3333  # J98: 4000 -> 3933.3333
$ws.Cells.Item(98,12).Value = 3933.3333  # L98: 4000 -> 3933.3333
$ws.Cells.Item(98,14).Value = -6929.3333  # N98: -6996 -> -6929.3333
$ws.Cells.Item(122,8).Value = 90911980  # H122: 111113784 -> 90911980
$ws.Cells.Item(122,10).Value = 3933.3333  # J122: 4000 -> 3933.3333
$ws.Cells.Item(122,12).Value = 11799.9999  # L122: 12000 -> 11799.9999
$ws.Cells.Item(122,14).Value = -16699.9999  # N122: -16900 -> -16699.9999
$ws.Cells.Item(137,8).Value = 5552.3076  # H137: 6610.727 -> 5552.3076
$ws.Cells.Item(137,9).Value = 3834.7273  # I137: 5152.5 -> 3834.7273
$ws.Cells.Item(137,10).Value = 14999  # J137: 10499.333 -> 14999
$ws.Cells.Item(137,11).Value = 11504.1819  # K137: 15457.5 -> 11504.1819
$ws.Cells.Item(137,12).Value = 44997  # L137: 31497.999 -> 44997
$ws.Cells.Item(137,13).Value = -8954.1819  # M137: -12907.5 -> -8954.1819
$ws.Cells.Item(137,14).Value = -50097  # N137: -36597.999 -> -50097

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32,8).Value = 8479756  # H32: 8777289 -> 8479756
$ws.Cells.Item(32,9).Value = 11365684  # I32: 11906904 -> 11365684
$ws.Cells.Item(32,11).Value = 11365684  # K32: 11906904 -> 11365684
$ws.Cells.Item(32,13).Value = -11365397  # M32: -11906617 -> -11365397
$ws.Cells.Item(45,8).Value = 33336070  # H45: 25002392 -> 33336070
$ws.Cells.Item(45,9).Value = 50002010  # I45: 33335124 -> 50002010
$ws.Cells.Item(45,11).Value = 50002010  # K45: 33335124 -> 50002010
$ws.Cells.Item(45,13).Value = -50001633  # M45: -33334747 -> -50001633
$ws.Cells.Item(61,8).Value = 39482052  # H61: 44126564 -> 39482052
$ws.Cells.Item(61,9).Value = 45461184  # I61: 62507750 -> 45461184
$ws.Cells.Item(61,10).Value = 31260742  # J61: 27787726 -> 31260742
$ws.Cells.Item(61,11).Value = 45461184  # K61: 62507750 -> 45461184
$ws.Cells.Item(61,12).Value = 31260742  # L61: 27787726 -> 31260742
$ws.Cells.Item(61,13).Value = -45460972  # M61: -62507538 -> -45460972
$ws.Cells.Item(61,14).Value = -31261166  # N61: -27788150 -> -31261166
$ws.Cells.Item(63,8).Value = 4635.2383  # H63: 4654.1904 -> 4635.2383
$ws.Cells.Item(63,9).Value = 2256.077  # I63: 2286.6924 -> 2256.077
$ws.Cells.Item(63,11).Value = 2256.077  # K63: 2286.6924 -> 2256.077
$ws.Cells.Item(63,13).Value = -1570.077  # M63: -1600.6924 -> -1570.077
$ws.Cells.Item(66,8).Value = 4635.2383  # H66: 4654.1904 -> 4635.2383
$ws.Cells.Item(66,9).Value = 2256.077  # I66: 2286.6924 -> 2256.077
$ws.Cells.Item(66,11).Value = 11280.385  # K66: 11433.462 -> 11280.385
$ws.Cells.Item(66,13).Value = -7848.385000000002  # M66: -8001.462 -> -7848.385000000002
$ws.Cells.Item(92,8).Value = 51148.168  # H92: 51997.8 -> 51148.168
$ws.Cells.Item(92,10).Value = 54377.8  # J92: 56247.25 -> 54377.8
$ws.Cells.Item(92,12).Value = 54377.8  # L92: 56247.25 -> 54377.8
$ws.Cells.Item(92,14).Value = -59369.8  # N92: -61239.25 -> -59369.8
$ws.Cells.Item(122,8).Value = 3315.2144  # H122: 4019.3635 -> 3315.2144
$ws.Cells.Item(122,9).Value = 2379.2222  # I122: 3202.1667 -> 2379.2222
$ws.Cells.Item(122,11).Value = 7137.6666  # K122: 9606.500100000001 -> 7137.6666
$ws.Cells.Item(122,13).Value = -4687.6666  # M122: -7156.500100000001 -> -4687.6666
$ws.Cells.Item(133,8).Value = 69998  # H133: 66998.625 -> 69998
$ws.Cells.Item(133,9).Value = 0  # I133: 69997 -> 0
$ws.Cells.Item(133,10).Value = 69998  # J133: 66570.28999999999 -> 69998
$ws.Cells.Item(133,11).Value = 0  # K133: 69997 -> 0
$ws.Cells.Item(133,12).Value = 69998  # L133: 66570.28999999999 -> 69998
$ws.Cells.Item(133,13).ClearContents()  # M133: remove (was -67467)
$ws.Cells.Item(133,14).Value = -75058  # N133: -71630.28999999999 -> -75058
$ws.Cells.Item(136,8).Value = 39482052  # H136: 44126564 -> 39482052
$ws.Cells.Item(136,9).Value = 45461184  # I136: 62507750 -> 45461184
$ws.Cells.Item(136,10).Value = 31260742  # J136: 27787726 -> 31260742
$ws.Cells.Item(136,11).Value = 136383552  # K136: 187523250 -> 136383552
$ws.Cells.Item(136,12).Value = 93782226  # L136: 83363178 -> 93782226
$ws.Cells.Item(136,13).Value = -136381002  # M136: -187520700 -> -136381002
$ws.Cells.Item(136,14).Value = -93787326  # N136: -83368278 -> -93787326

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(10,8).Value = 0  # H10: 5 -> 0
$ws.Cells.Item(10,9).Value = 0  # I10: 5 -> 0
$ws.Cells.Item(10,11).Value = 0  # K10: 5 -> 0
$ws.Cells.Item(10,13).ClearContents()  # M10: remove (was 135)
$ws.Cells.Item(86,8).Value = 2424.9656  # H86: 2676.7144 -> 2424.9656
$ws.Cells.Item(86,9).Value = 2704.1904  # I86: 2854.3635 -> 2704.1904
$ws.Cells.Item(86,10).Value = 1692  # J86: 2025.3334 -> 1692
$ws.Cells.Item(86,11).Value = 2704.1904  # K86: 2854.3635 -> 2704.1904
$ws.Cells.Item(86,12).Value = 1692  # L86: 2025.3334 -> 1692
$ws.Cells.Item(86,13).Value = -1581.1904  # M86: -1731.3635 -> -1581.1904
$ws.Cells.Item(86,14).Value = -3938  # N86: -4271.3334 -> -3938
$ws.Cells.Item(89,8).Value = 2424.9656  # H89: 2676.7144 -> 2424.9656
$ws.Cells.Item(89,9).Value = 2704.1904  # I89: 2854.3635 -> 2704.1904
$ws.Cells.Item(89,10).Value = 1692  # J89: 2025.3334 -> 1692
$ws.Cells.Item(89,11).Value = 13520.952  # K89: 14271.8175 -> 13520.952
$ws.Cells.Item(89,12).Value = 8460  # L89: 10126.667 -> 8460
$ws.Cells.Item(89,13).Value = -7904.951999999999  # M89: -8655.817499999999 -> -7904.951999999999
$ws.Cells.Item(89,14).Value = -19692  # N89: -21358.667 -> -19692

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16,8).Value = 820.3333  # H16: 841.9091 -> 820.3333
$ws.Cells.Item(16,9).Value = 820.3333  # I16: 841.9091 -> 820.3333
$ws.Cells.Item(16,11).Value = 820.3333  # K16: 841.9091 -> 820.3333
$ws.Cells.Item(16,13).Value = -533.3333  # M16: -554.9091 -> -533.3333
$ws.Cells.Item(31,8).Value = 1460784.2  # H31: 1460672.4 -> 1460784.2
$ws.Cells.Item(31,9).Value = 3459.2  # I31: 3384.5 -> 3459.2
$ws.Cells.Item(31,10).Value = 2123204.8  # J31: 2335045 -> 2123204.8
$ws.Cells.Item(31,11).Value = 3459.2  # K31: 3384.5 -> 3459.2
$ws.Cells.Item(31,12).Value = 2123204.8  # L31: 2335045 -> 2123204.8
$ws.Cells.Item(31,13).Value = -3164.2  # M31: -3089.5 -> -3164.2
$ws.Cells.Item(31,14).Value = -2123794.8  # N31: -2335635 -> -2123794.8
$ws.Cells.Item(34,8).Value = 1460784.2  # H34: 1460672.4 -> 1460784.2
$ws.Cells.Item(34,9).Value = 3459.2  # I34: 3384.5 -> 3459.2
$ws.Cells.Item(34,10).Value = 2123204.8  # J34: 2335045 -> 2123204.8
$ws.Cells.Item(34,11).Value = 3459.2  # K34: 3384.5 -> 3459.2
$ws.Cells.Item(34,12).Value = 2123204.8  # L34: 2335045 -> 2123204.8
$ws.Cells.Item(34,13).Value = -3257.2  # M34: -3182.5 -> -3257.2
$ws.Cells.Item(34,14).Value = -2123608.8  # N34: -2335449 -> -2123608.8
$ws.Cells.Item(58,8).Value = 9130.571  # H58: 8613.625 -> 9130.571
$ws.Cells.Item(58,10).Value = 19000  # J58: 11997.5 -> 19000
$ws.Cells.Item(58,12).Value = 19000  # L58: 11997.5 -> 19000
$ws.Cells.Item(58,14).Value = -19406  # N58: -12403.5 -> -19406
$ws.Cells.Item(113,8).Value = 820.3333  # H113: 841.9091 -> 820.3333
$ws.Cells.Item(113,9).Value = 820.3333  # I113: 841.9091 -> 820.3333
$ws.Cells.Item(113,11).Value = 820.3333  # K113: 841.9091 -> 820.3333
$ws.Cells.Item(113,13).Value = 1349.6667  # M113: 1328.0909 -> 1349.6667
$ws.Cells.Item(132,8).Value = 11317.25  # H132: 10481.777 -> 11317.25
$ws.Cells.Item(132,10).Value = 13321.5  # J132: 11961 -> 13321.5
$ws.Cells.Item(132,12).Value = 39964.5  # L132: 35883 -> 39964.5
$ws.Cells.Item(132,14).Value = -45024.5  # N132: -40943 -> -45024.5
$ws.Cells.Item(134,8).Value = 4560.625  # H134: 4386 -> 4560.625
$ws.Cells.Item(134,9).Value = 3747.5  # I134: 3639.1428 -> 3747.5
$ws.Cells.Item(134,11).Value = 11242.5  # K134: 10917.4284 -> 11242.5
$ws.Cells.Item(134,13).Value = -8707.5  # M134: -8382.428400000001 -> -8707.5
$ws.Cells.Item(136,8).Value = 9130.571  # H136: 8613.625 -> 9130.571
$ws.Cells.Item(136,10).Value = 19000  # J136: 11997.5 -> 19000
$ws.Cells.Item(136,12).Value = 57000  # L136: 35992.5 -> 57000
$ws.Cells.Item(136,14).Value = -62100  # N136: -41092.5 -> -62100

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(37,8).Value = 58166  # H37: 53624.25 -> 58166
$ws.Cells.Item(37,10).Value = 58166  # J37: 53624.25 -> 58166
$ws.Cells.Item(37,12).Value = 174498  # L37: 160872.75 -> 174498
$ws.Cells.Item(37,14).Value = -174722  # N37: -161096.75 -> -174722
$ws.Cells.Item(38,8).Value = 103.47369  # H38: 111.117645 -> 103.47369
$ws.Cells.Item(38,9).Value = 109.583336  # I38: 116.90909 -> 109.583336
$ws.Cells.Item(38,10).Value = 93  # J38: 100.5 -> 93
$ws.Cells.Item(38,11).Value = 328.750008  # K38: 350.72727 -> 328.750008
$ws.Cells.Item(38,12).Value = 279  # L38: 301.5 -> 279
$ws.Cells.Item(38,13).Value = 18.24999200000002  # M38: -3.727270000000033 -> 18.24999200000002
$ws.Cells.Item(38,14).Value = -973  # N38: -995.5 -> -973
$ws.Cells.Item(121,8).Value = 1232.4138  # H121: 819.7241 -> 1232.4138
$ws.Cells.Item(121,10).Value = 1811.9375  # J121: 1063.9375 -> 1811.9375
$ws.Cells.Item(121,12).Value = 5435.8125  # L121: 3191.8125 -> 5435.8125
$ws.Cells.Item(121,14).Value = -8055.8125  # N121: -5811.8125 -> -8055.8125

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70,8).Value = 20874.5  # H70: 20874.625 -> 20874.5
$ws.Cells.Item(70,9).Value = 28199.2  # I70: 28199.4 -> 28199.2
$ws.Cells.Item(70,11).Value = 28199.2  # K70: 28199.4 -> 28199.2
$ws.Cells.Item(70,13).Value = -27929.2  # M70: -27929.4 -> -27929.2
$ws.Cells.Item(73,8).Value = 20874.5  # H73: 20874.625 -> 20874.5
$ws.Cells.Item(73,9).Value = 28199.2  # I73: 28199.4 -> 28199.2
$ws.Cells.Item(73,11).Value = 28199.2  # K73: 28199.4 -> 28199.2
$ws.Cells.Item(73,13).Value = -27263.2  # M73: -27263.4 -> -27263.2
$ws.Cells.Item(102,8).Value = 2817.7036  # H102: 2963.28 -> 2817.7036
$ws.Cells.Item(102,9).Value = 2154.6  # I102: 2283.111 -> 2154.6
$ws.Cells.Item(102,11).Value = 2154.6  # K102: 2283.111 -> 2154.6
$ws.Cells.Item(102,13).Value = -532.5999999999999  # M102: -661.1109999999999 -> -532.5999999999999
$ws.Cells.Item(113,8).Value = 3911.5833  # H113: 4089.6191 -> 3911.5833
$ws.Cells.Item(113,9).Value = 2823.5  # I113: 2941.7144 -> 2823.5
$ws.Cells.Item(113,10).Value = 4455.625  # J113: 4663.5713 -> 4455.625
$ws.Cells.Item(113,11).Value = 2823.5  # K113: 2941.7144 -> 2823.5
$ws.Cells.Item(113,12).Value = 4455.625  # L113: 4663.5713 -> 4455.625
$ws.Cells.Item(113,13).Value = -653.5  # M113: -771.7143999999998 -> -653.5
$ws.Cells.Item(113,14).Value = -8795.625  # N113: -9003.5713 -> -8795.625
$ws.Cells.Item(132,8).Value = 66674976  # H132: 90919910 -> 66674976
$ws.Cells.Item(132,9).Value = 90910780  # I132: 125001816 -> 90910780
$ws.Cells.Item(132,10).Value = 26502.75  # J132: 34837.668 -> 26502.75
$ws.Cells.Item(132,11).Value = 272732340  # K132: 375005448 -> 272732340
$ws.Cells.Item(132,12).Value = 79508.25  # L132: 104513.004 -> 79508.25
$ws.Cells.Item(132,13).Value = -272729810  # M132: -375002918 -> -272729810
$ws.Cells.Item(132,14).Value = -84568.25  # N132: -109573.004 -> -84568.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22,8).Value = 1772.9412  # H22: 1943.4667 -> 1772.9412
$ws.Cells.Item(22,9).Value = 1918.3  # I22: 2094.4443 -> 1918.3
$ws.Cells.Item(22,10).Value = 1565.2858  # J22: 1717 -> 1565.2858
$ws.Cells.Item(22,11).Value = 1918.3  # K22: 2094.4443 -> 1918.3
$ws.Cells.Item(22,12).Value = 1565.2858  # L22: 1717 -> 1565.2858
$ws.Cells.Item(22,13).Value = -1623.3  # M22: -1799.4443 -> -1623.3
$ws.Cells.Item(22,14).Value = -2155.2858  # N22: -2307 -> -2155.2858
$ws.Cells.Item(27,8).Value = 1772.9412  # H27: 1943.4667 -> 1772.9412
$ws.Cells.Item(27,9).Value = 1918.3  # I27: 2094.4443 -> 1918.3
$ws.Cells.Item(27,10).Value = 1565.2858  # J27: 1717 -> 1565.2858
$ws.Cells.Item(27,11).Value = 1918.3  # K27: 2094.4443 -> 1918.3
$ws.Cells.Item(27,12).Value = 1565.2858  # L27: 1717 -> 1565.2858
$ws.Cells.Item(27,13).Value = -1811.3  # M27: -1987.4443 -> -1811.3
$ws.Cells.Item(27,14).Value = -1779.2858  # N27: -1931 -> -1779.2858
$ws.Cells.Item(93,8).Value = 47620372  # H93: 47620376 -> 47620372
$ws.Cells.Item(93,10).Value = 1650.7273  # J93: 1651.8182 -> 1650.7273
$ws.Cells.Item(93,12).Value = 1650.7273  # L93: 1651.8182 -> 1650.7273
$ws.Cells.Item(93,14).Value = -4146.7273  # N93: -4147.8182 -> -4146.7273
$ws.Cells.Item(122,8).Value = 6231.2324  # H122: 6032.3477 -> 6231.2324
$ws.Cells.Item(122,9).Value = 5790.758  # I122: 6069.1934 -> 5790.758
$ws.Cells.Item(122,10).Value = 7684.8  # J122: 5956.2 -> 7684.8
$ws.Cells.Item(122,11).Value = 17372.274  # K122: 18207.5802 -> 17372.274
$ws.Cells.Item(122,12).Value = 23054.4  # L122: 17868.6 -> 23054.4
$ws.Cells.Item(122,13).Value = -14922.274  # M122: -15757.5802 -> -14922.274
$ws.Cells.Item(122,14).Value = -27954.4  # N122: -22768.6 -> -27954.4
$ws.Cells.Item(133,8).Value = 68849.5  # H133: 68566.336 -> 68849.5
$ws.Cells.Item(133,10).Value = 68849.5  # J133: 68566.336 -> 68849.5
$ws.Cells.Item(133,12).Value = 68849.5  # L133: 68566.336 -> 68849.5
$ws.Cells.Item(133,14).Value = -73909.5  # N133: -73626.336 -> -73909.5
$ws.Cells.Item(136,8).Value = 167000.56  # H136: 107814.36 -> 167000.56
$ws.Cells.Item(136,9).Value = 31250  # I136: 16049.5 -> 31250
$ws.Cells.Item(136,10).Value = 275601  # J136: 230167.5 -> 275601
$ws.Cells.Item(136,11).Value = 93750  # K136: 48148.5 -> 93750
$ws.Cells.Item(136,12).Value = 826803  # L136: 690502.5 -> 826803
$ws.Cells.Item(136,13).Value = -91200  # M136: -45598.5 -> -91200
$ws.Cells.Item(136,14).Value = -831903  # N136: -695602.5 -> -831903

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(82,8).Value = 34746.75  # H82: 33572 -> 34746.75
$ws.Cells.Item(82,9).Value = 28990  # I82: 28931.5 -> 28990
$ws.Cells.Item(82,11).Value = 28990  # K82: 28931.5 -> 28990
$ws.Cells.Item(82,13).Value = -28607  # M82: -28548.5 -> -28607
$ws.Cells.Item(85,8).Value = 34746.75  # H85: 33572 -> 34746.75
$ws.Cells.Item(85,9).Value = 28990  # I85: 28931.5 -> 28990
$ws.Cells.Item(85,11).Value = 28990  # K85: 28931.5 -> 28990
$ws.Cells.Item(85,13).Value = -27664  # M85: -27605.5 -> -27664
$ws.Cells.Item(88,8).Value = 10000  # H88: 17971 -> 10000
$ws.Cells.Item(88,9).Value = 0  # I88: 17971 -> 0
$ws.Cells.Item(88,10).Value = 10000  # J88: 0 -> 10000
$ws.Cells.Item(88,11).Value = 0  # K88: 17971 -> 0
$ws.Cells.Item(88,12).Value = 10000  # L88: 0 -> 10000
$ws.Cells.Item(88,13).ClearContents()  # M88: remove (was -17565)
$ws.Cells.Item(88,14).Value = -10812  # N88: None -> -10812
$ws.Cells.Item(91,8).Value = 10000  # H91: 17971 -> 10000
$ws.Cells.Item(91,9).Value = 0  # I91: 17971 -> 0
$ws.Cells.Item(91,10).Value = 10000  # J91: 0 -> 10000
$ws.Cells.Item(91,11).Value = 0  # K91: 17971 -> 0
$ws.Cells.Item(91,12).Value = 10000  # L91: 0 -> 10000
$ws.Cells.Item(91,13).ClearContents()  # M91: remove (was -16567)
$ws.Cells.Item(91,14).Value = -12808  # N91: None -> -12808
